{"js": "// \"... en faisant preuve d'esprit d'\u00e9quipe.\" -> bold \"esprit d'\u00e9quipe\"\n// (the leading \"d'\" right before \"esprit\" stays un-bolded).\n//\n// Word pairs <w:b/> with <w:bCs/> whenever bold is applied (this document's\n// own pre-existing bold run \"Mais\" already has both). The host's\n// font.boldBidirectional setter (the Office.js analogue of the COM\n// Font.BoldBi flag) is unfortunately bugged in this runtime: instead of\n// scoping to the target range it stamps <w:bCs/> onto every run in the\n// whole document body. So we can't use it here and instead inject the\n// exact run/rPr XML we want for just the matched span.\n\nconst body = context.document.body;\n\nconst results = body.search(\"esprit d\\u2019\\u00e9quipe\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find \\u201cesprit d\\u2019\\u00e9quipe\\u201d in the document.\");\n}\n\nconst target = results.items[0];\n\n// Rebuild just this span as three runs (mirroring the existing run\n// boundaries: \"esprit\" | \" \" | \"d'\u00e9quipe\", separated by the proofErr marks\n// already present around \"d'\u00e9quipe\"), each bold with both <w:b/> and\n// <w:bCs/>.\nconst innerRuns =\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>esprit</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>d\\u2019\\u00e9quipe</w:t></w:r>';\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' + innerRuns + '</w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the phrase \"esprit d'\u00e9quipe\" (curly apostrophe, as used in the\n# document) within \"... en faisant preuve d'esprit d'\u00e9quipe.\" and make it\n# bold. The leading \"d'\" (before \"esprit\") stays un-bolded.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"esprit d\u2019\u00e9quipe\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    $rng.Font.Bold = 1\n}\n"}
